$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "601.70", "8.00") must be
# forced to the Text number format first, otherwise Excel auto-converts the string
# into a floating point value (losing formatting like trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D5").Value = "601.70"
$ws.Range("D6").Value = "178.42"
$ws.Range("D8").Value = "0.524"
$ws.Range("D10").Value = "0.173"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("D17").Value = "26.61"
$ws.Range("D19").Value = "11.96"
$ws.Range("D20").Value = "8.00"
$ws.Range("D21").Value = "379.94"
$ws.Range("D22").Value = "4.20"
$ws.Range("D23").Value = "2.05"
$ws.Range("D24").Value = "72.47"
$ws.Range("D26").Value = "4.37"
$ws.Range("D27").Value = "9.96"
$ws.Range("D31").Value = "8.17"
$ws.Range("D32").Value = "520.73"
$ws.Range("D35").Value = "0.999"
$ws.Range("D36").Value = "164.75"
$ws.Range("D37").Value = "19.54"
$ws.Range("D38").Value = "19.11"
$ws.Range("D41").Value = "1.84"
$ws.Range("D43").Value = "5.06"
$ws.Range("D44").Value = "2.59"
$ws.Range("D45").Value = "0.334"
$ws.Range("D46").Value = "39.31"
$ws.Range("D47").Value = "152.93"
$ws.Range("D48").Value = "3.75"
$ws.Range("D49").Value = "0.550"

# Remaining cells (coin names, links, prices that already contain non-numeric
# characters such as extra "." separators or percentage text, and the volume
# column) can be assigned directly -- Excel keeps them as text automatically.
$ws.Range("D2").Value = "72.151.34"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "2.663.78"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.661.66"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("E10").Value = "  +5.21%  "
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "3.151.01"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").Value = "72.002.19"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "2.668.26"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("E20").Value = "  +3.16%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("E23").Value = "  +11.09%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").Value = "2.809.57"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "0.0₃0946"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  -6.84%  "
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("E51").Value = "  +1.53%  "
